$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1433.9
$ws.Range("J33").Value = 232.66667
$ws.Range("L33").Value = 232.66667
$ws.Range("N33").Value = -690.6666700000001

$ws.Range("H45").Value = 5000
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -15384

$ws.Range("H112").Value = 3038.2222
$ws.Range("I112").Value = 840
$ws.Range("J112").Value = 3883.6924
$ws.Range("K112").Value = 2520
$ws.Range("L112").Value = 11651.0772
$ws.Range("M112").Value = -1412
$ws.Range("N112").Value = -13867.0772

$ws.Range("H121").Value = 1771.1666
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1771.1666
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 5313.4998
$ws.Range("N121").Value = -8807.4998
$ws.Range("M121").ClearContents()

$ws.Range("H138").Value = 3740.0876
$ws.Range("I138").Value = 1309.4348
$ws.Range("J138").Value = 5384.353
$ws.Range("K138").Value = 3928.3044
$ws.Range("L138").Value = 16153.059
$ws.Range("M138").Value = 1211.6956
$ws.Range("N138").Value = -26433.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2400.8667
$ws.Range("I2").Value = 1176.25
$ws.Range("J2").Value = 3800.4285
$ws.Range("K2").Value = 1176.25
$ws.Range("L2").Value = 3800.4285
$ws.Range("M2").Value = -1063.25
$ws.Range("N2").Value = -4026.4285

$ws.Range("H32").Value = 3977756.5
$ws.Range("I32").Value = 5288.971
$ws.Range("J32").Value = 22251106
$ws.Range("K32").Value = 5288.971
$ws.Range("L32").Value = 22251106
$ws.Range("M32").Value = -5001.971
$ws.Range("N32").Value = -22251680

$ws.Range("H61").Value = 4936.8184
$ws.Range("I61").Value = 2559.4285
$ws.Range("J61").Value = 9097.25
$ws.Range("K61").Value = 2559.4285
$ws.Range("L61").Value = 9097.25
$ws.Range("M61").Value = -2347.4285
$ws.Range("N61").Value = -9521.25

$ws.Range("H116").Value = 2400.8667
$ws.Range("I116").Value = 1176.25
$ws.Range("J116").Value = 3800.4285
$ws.Range("K116").Value = 1176.25
$ws.Range("L116").Value = 3800.4285
$ws.Range("M116").Value = 1117.75
$ws.Range("N116").Value = -8388.4285

$ws.Range("H136").Value = 4936.8184
$ws.Range("I136").Value = 2559.4285
$ws.Range("J136").Value = 9097.25
$ws.Range("K136").Value = 7678.2855
$ws.Range("L136").Value = 27291.75
$ws.Range("M136").Value = -5128.2855
$ws.Range("N136").Value = -32391.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2400.8667
$ws.Range("I3").Value = 1176.25
$ws.Range("J3").Value = 3800.4285
$ws.Range("K3").Value = 1176.25
$ws.Range("L3").Value = 3800.4285
$ws.Range("M3").Value = -1062.25
$ws.Range("N3").Value = -4028.4285

$ws.Range("H7").Value = 499.63635
$ws.Range("I7").Value = 156.57143
$ws.Range("J7").Value = 1100
$ws.Range("K7").Value = 156.57143
$ws.Range("L7").Value = 1100
$ws.Range("M7").Value = -43.57142999999999
$ws.Range("N7").Value = -1326

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H99").Value = 4264.636
$ws.Range("I99").Value = 3900
$ws.Range("J99").Value = 4902.75
$ws.Range("K99").Value = 3900
$ws.Range("L99").Value = 4902.75
$ws.Range("M99").Value = -2402
$ws.Range("N99").Value = -7898.75

$ws.Range("H134").Value = 6765.5864
$ws.Range("I134").Value = 3927
$ws.Range("J134").Value = 7846.952
$ws.Range("K134").Value = 11781
$ws.Range("L134").Value = 23540.856
$ws.Range("M134").Value = -9246
$ws.Range("N134").Value = -28610.856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1300
$ws.Range("I10").Value = 1300
$ws.Range("K10").Value = 1300
$ws.Range("M10").Value = -1161

$ws.Range("H68").Value = 32950
$ws.Range("J68").Value = 32950
$ws.Range("L68").Value = 32950
$ws.Range("N68").Value = -34448

$ws.Range("H71").Value = 32950
$ws.Range("J71").Value = 32950
$ws.Range("L71").Value = 98850
$ws.Range("N71").Value = -106338

$ws.Range("H119").Value = 49400
$ws.Range("J119").Value = 49400
$ws.Range("L119").Value = 49400
$ws.Range("N119").Value = -59076

$ws.Range("H132").Value = 6670524.5
$ws.Range("I132").Value = 11113430
$ws.Range("J132").Value = 6166.3335
$ws.Range("K132").Value = 33340290
$ws.Range("L132").Value = 18499.0005
$ws.Range("M132").Value = -33337760
$ws.Range("N132").Value = -23559.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1184.25
$ws.Range("I13").Value = 1011
$ws.Range("J13").Value = 1308
$ws.Range("K13").Value = 1011
$ws.Range("L13").Value = 1308
$ws.Range("M13").Value = -872
$ws.Range("N13").Value = -1586

$ws.Range("H132").Value = 2187.6667
$ws.Range("I132").Value = 1906.625
$ws.Range("J132").Value = 2749.75
$ws.Range("K132").Value = 5719.875
$ws.Range("L132").Value = 8249.25
$ws.Range("M132").Value = -3189.875
$ws.Range("N132").Value = -13309.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 66672120
$ws.Range("I132").Value = 111114530
$ws.Range("J132").Value = 8499.666999999999
$ws.Range("K132").Value = 333343590
$ws.Range("L132").Value = 25499.001
$ws.Range("M132").Value = -333341060
$ws.Range("N132").Value = -30559.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 2777.25
$ws.Range("I132").Value = 1915.2778
$ws.Range("J132").Value = 4328.8
$ws.Range("K132").Value = 5745.8334
$ws.Range("L132").Value = 12986.4
$ws.Range("M132").Value = -3215.8334
$ws.Range("N132").Value = -18046.4

$ws.Range("H136").Value = 5682675
$ws.Range("I136").Value = 6250802
$ws.Range("J136").Value = 1404.5
$ws.Range("K136").Value = 18752406
$ws.Range("L136").Value = 4213.5
$ws.Range("M136").Value = -18749856
$ws.Range("N136").Value = -9313.5
